$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 3 and row 4 for columns D, J, K, L, M, P
$cols = @("D", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $cell3 = $ws.Range($col + "3")
    $cell4 = $ws.Range($col + "4")
    $v3 = $cell3.Value2
    $v4 = $cell4.Value2
    $cell3.Value2 = $v4
    $cell4.Value2 = $v3
}
